# Update the division-problem answers table.
# The document contains a single 5-column table; rows 1, 5, 9, 13, 17
# (1-indexed) hold the visible "a÷b=c, d" strings (the intervening rows
# are blank spacer rows). Several of the original strings repeat
# (e.g. "62÷8=7, 6" and "81÷6=13, 3" each appear twice), so addressing
# cells directly via Table.Cell(row, col) is used instead of a global
# Find/Replace, which would not be able to distinguish the duplicates.
#
# Assigning to Cell.Range.Text replaces only the text of the existing
# run, preserving its rPr (rFonts/sz) formatting, matching the diff.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "98÷9=10, 8"
$t.Cell(1, 2).Range.Text  = "72÷7=10, 2"
$t.Cell(1, 3).Range.Text  = "59÷2=29, 1"
$t.Cell(1, 4).Range.Text  = "33÷4=8, 1"
$t.Cell(1, 5).Range.Text  = "82÷8=10, 2"

$t.Cell(5, 1).Range.Text  = "41÷5=8, 1"
$t.Cell(5, 2).Range.Text  = "90÷6=15, 0"
$t.Cell(5, 3).Range.Text  = "99÷5=19, 4"
$t.Cell(5, 4).Range.Text  = "74÷7=10, 4"
$t.Cell(5, 5).Range.Text  = "17÷2=8, 1"

$t.Cell(9, 1).Range.Text  = "78÷5=15, 3"
$t.Cell(9, 2).Range.Text  = "70÷8=8, 6"
$t.Cell(9, 3).Range.Text  = "44÷2=22, 0"
$t.Cell(9, 4).Range.Text  = "80÷6=13, 2"
$t.Cell(9, 5).Range.Text  = "18÷2=9, 0"

$t.Cell(13, 1).Range.Text = "57÷4=14, 1"
$t.Cell(13, 2).Range.Text = "54÷7=7, 5"
$t.Cell(13, 3).Range.Text = "13÷9=1, 4"
$t.Cell(13, 4).Range.Text = "96÷3=32, 0"
$t.Cell(13, 5).Range.Text = "25÷2=12, 1"

$t.Cell(17, 1).Range.Text = "59÷8=7, 3"
$t.Cell(17, 2).Range.Text = "65÷7=9, 2"
$t.Cell(17, 3).Range.Text = "68÷2=34, 0"
$t.Cell(17, 4).Range.Text = "43÷6=7, 1"
$t.Cell(17, 5).Range.Text = "95÷7=13, 4"
